# Update the "Fitness" column (C) values in Sheet1 to reflect a new run's
# logged results. The Run/Generation columns (A/B) are unchanged; only the
# Fitness numbers in column C (rows 2-252) are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13218
$ws.Range("C3").Value = 12612
$ws.Range("C4:C8").Value = 11105
$ws.Range("C9").Value = 10236
$ws.Range("C10:C12").Value = 8945
$ws.Range("C13:C33").Value = 7343
$ws.Range("C34:C252").Value = 7293
